$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-13 Friday" "2025-06-14 Saturday"
Replace-Text "67÷9=7, 4" "66÷5=13, 1"
Replace-Text "10÷3=3, 1" "90÷3=30, 0"
Replace-Text "86÷9=9, 5" "54÷7=7, 5"
Replace-Text "83÷2=41, 1" "87÷4=21, 3"
Replace-Text "27÷5=5, 2" "61÷5=12, 1"
Replace-Text "78÷7=11, 1" "46÷6=7, 4"
Replace-Text "43÷9=4, 7" "39÷6=6, 3"
Replace-Text "23÷7=3, 2" "50÷6=8, 2"
Replace-Text "52÷4=13, 0" "38÷8=4, 6"
Replace-Text "19÷7=2, 5" "66÷2=33, 0"
Replace-Text "26÷2=13, 0" "50÷8=6, 2"
Replace-Text "86÷5=17, 1" "63÷9=7, 0"
Replace-Text "79÷9=8, 7" "53÷7=7, 4"
Replace-Text "63÷6=10, 3" "51÷6=8, 3"
Replace-Text "94÷4=23, 2" "40÷3=13, 1"
Replace-Text "59÷2=29, 1" "10÷8=1, 2"
Replace-Text "75÷3=25, 0" "56÷6=9, 2"
Replace-Text "26÷4=6, 2" "76÷6=12, 4"
Replace-Text "77÷7=11, 0" "64÷5=12, 4"
Replace-Text "63÷3=21, 0" "54÷3=18, 0"
Replace-Text "99÷2=49, 1" "67÷5=13, 2"
Replace-Text "90÷6=15, 0" "49÷7=7, 0"
Replace-Text "83÷3=27, 2" "92÷6=15, 2"
Replace-Text "64÷7=9, 1" "93÷5=18, 3"
Replace-Text "36÷4=9, 0" "13÷3=4, 1"
